$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -7
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = -1
$ws.Range("F9").Value = 1
$ws.Range("F10").Value = 1
$ws.Range("F20").Value = 4
$ws.Range("F21").Value = 3
$ws.Range("F22").Value = -4
$ws.Range("F33").Value = 0
$ws.Range("F35").Value = 0
$ws.Range("F38").Value = -11
$ws.Range("F40").Value = -6
$ws.Range("F42").Value = -10
$ws.Range("F43").Value = -10
$ws.Range("F45").Value = -6
$ws.Range("F47").Value = -8
$ws.Range("F50").Value = 5
$ws.Range("F51").Value = -1
$ws.Range("F52").Value = -3
